# Insert a new weekly price-report row for "Haba" at row 68 (Femacal de La
# Calera, Coquimbo). This pushes the existing rows 68-140 down to 69-141
# (dimension grows from A1:R140 to A1:R141) and fills the freshly inserted
# row 68 with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("68").Insert()

$ws.Range("A68").Value = 3
$ws.Range("B68").Value = "Femacal de La Calera"
$ws.Range("C68").Value = "Coquimbo"
$ws.Range("D68").Value = 44721
$ws.Range("E68").Value = 5
$ws.Range("F68").Value = 100112026
$ws.Range("G68").Value = "Haba"
$ws.Range("H68").Value = "Sin especificar"
$ws.Range("I68").Value = "Primera"
$ws.Range("J68").Value = 76
$ws.Range("K68").Value = 21500
$ws.Range("L68").Value = 22000
$ws.Range("M68").Value = 21750
$ws.Range("N68").Value = "$/saco 25 kilos"
$ws.Range("O68").Value = "Provincia de Limarí"
$ws.Range("P68").Value = 870
$ws.Range("Q68").Value = 25
$ws.Range("R68").Value = "Hortaliza"
